$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "34÷3=11, 1"
$t.Cell(1, 2).Range.Text = "75÷2=37, 1"
$t.Cell(1, 3).Range.Text = "13÷2=6, 1"
$t.Cell(1, 4).Range.Text = "90÷6=15, 0"
$t.Cell(1, 5).Range.Text = "14÷4=3, 2"
$t.Cell(5, 1).Range.Text = "75÷7=10, 5"
$t.Cell(5, 2).Range.Text = "99÷3=33, 0"
$t.Cell(5, 3).Range.Text = "85÷7=12, 1"
$t.Cell(5, 4).Range.Text = "86÷3=28, 2"
$t.Cell(5, 5).Range.Text = "84÷5=16, 4"
$t.Cell(9, 1).Range.Text = "92÷5=18, 2"
$t.Cell(9, 2).Range.Text = "58÷8=7, 2"
$t.Cell(9, 3).Range.Text = "85÷7=12, 1"
$t.Cell(9, 4).Range.Text = "45÷9=5, 0"
$t.Cell(9, 5).Range.Text = "15÷2=7, 1"
$t.Cell(13, 1).Range.Text = "28÷4=7, 0"
$t.Cell(13, 2).Range.Text = "68÷9=7, 5"
$t.Cell(13, 3).Range.Text = "20÷2=10, 0"
$t.Cell(13, 4).Range.Text = "46÷7=6, 4"
$t.Cell(13, 5).Range.Text = "10÷3=3, 1"
$t.Cell(17, 1).Range.Text = "57÷6=9, 3"
$t.Cell(17, 2).Range.Text = "45÷8=5, 5"
$t.Cell(17, 3).Range.Text = "99÷7=14, 1"
$t.Cell(17, 4).Range.Text = "51÷6=8, 3"
$t.Cell(17, 5).Range.Text = "16÷3=5, 1"
